$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$sub = $tr.Characters(37, 1)
$sub.Font.Size = 25
